$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header (H1) — reuse the existing header formatting
# (bold, bordered, centered) from the adjacent "sum" header cell (G1).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save values for each data row (H2:H8)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("H8").Value = 0
